$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.832000000000001
$ws.Range("B10").Value = 6.039
$ws.Range("B12").Value = 4.951000000000001
$ws.Range("B18").Value = 5.137
$ws.Range("B25").Value = 5.510999999999999
$ws.Range("B37").Value = 8.657
$ws.Range("B55").Value = 4.763
$ws.Range("B68").Value = 5.138999999999999
$ws.Range("B77").Value = 5.459999999999999
$ws.Range("B78").Value = 7.410000000000001
$ws.Range("B79").Value = 5.040000000000001
$ws.Range("B80").Value = 7.961000000000001
$ws.Range("B81").Value = 5.799
$ws.Range("B82").Value = 5.457
$ws.Range("B84").Value = 5.82
$ws.Range("B101").Value = 5.470000000000001
$ws.Range("B102").Value = 7.211999999999999
